$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.502985000610352
$ws.Range("B1").Value = 2.737139701843262
$ws.Range("C1").Value = 1.806487083435059
$ws.Range("D1").Value = 1.568248152732849
$ws.Range("E1").Value = 1.538068294525146
